$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.970.59'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '1.654.10'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.17'
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0615'
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.18'
$ws.Range("E10").Value = '  +4.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0875'
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("D12").Value = '1.888.91'
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("D13").Value = '1.655.55'
$ws.Range("E13").Value = '  +2.98%  '
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").Value = '  +2.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.11'
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("D17").Value = '26.996.88'
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '235.62'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.71'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  +3.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.31'
$ws.Range("E23").Value = '  +3.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  +3.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.47'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.09'
$ws.Range("E26").Value = '  +1.57%  '
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.80'
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").Value = '1.549.25'
$ws.Range("E32").Value = '  +3.85%  '
$ws.Range("E33").Value = '  +2.80%  '
$ws.Range("E34").Value = '  +4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  +7.82%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.582'
$ws.Range("E37").Value = '  +3.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.892'
$ws.Range("E38").Value = '  +8.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0169'
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.99'
$ws.Range("E40").Value = '  +3.15%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  +2.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.42'
$ws.Range("E43").Value = '  +7.31%  '
$ws.Range("D44").Value = '1.795.52'
$ws.Range("E44").Value = '  +2.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.775'
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.917'
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.19'
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.51'
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  +4.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0984'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0505'
$ws.Range("E51").Value = '  +0.78%  '
